# Remove the trailing "Ver no Jupiter..." and "(c) 2020 ..." footer
# paragraphs (and the blank paragraph that separates them from the
# "Requisitos" section above), while leaving the blank paragraph and the
# page-break paragraph that follow them untouched.

$d = $word.ActiveDocument

$startText = "Ver no Jupiter Salvar em pdf Salvar em docx"
$endText   = "Powered by Jekyll and Github pages"

$count = $d.Paragraphs.Count
$startIndex = -1
$endIndex = -1

for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*$startText*") {
        $startIndex = $i
    }
    if ($t -like "*$endText*") {
        $endIndex = $i
    }
}

if ($startIndex -gt 0 -and $endIndex -ge $startIndex) {
    # Also swallow the blank "Normal" paragraph immediately preceding the
    # "Ver no Jupiter..." paragraph, so it disappears along with the footer.
    $deleteFrom = $startIndex
    $prev = $d.Paragraphs.Item($startIndex - 1)
    if ($prev.Range.Text.Trim() -eq "") {
        $deleteFrom = $startIndex - 1
    }

    $rangeStart = $d.Paragraphs.Item($deleteFrom).Range.Start
    $rangeEnd = $d.Paragraphs.Item($endIndex).Range.End
    $r = $d.Range($rangeStart, $rangeEnd)
    $r.Delete()
}
